# [Kadastro App] Yeni kayit eklendi: 2970
# Adds a new record row (row 38) to both the "Kayitlar" (master list) sheet
# and the "Erdemli" (district) sheet, mirroring the existing rows' layout.

$wb = $excel.ActiveWorkbook

$sheetNames = @("Kayitlar", "Erdemli")

$kayitNo    = "2970"
$tarih      = "2025-09-10"
$birim      = "Erdemli"
$parselSayi = "1"
$is         = "ÇAP"
$personel   = "AYHAN KARADAYI (K.Teknisyeni)"

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $row = 38

    # Leading apostrophe forces these (otherwise numeric/date-looking)
    # values to be stored as text, matching the rest of the column.
    $ws.Range("A$row").Value = "'" + $kayitNo
    $ws.Range("B$row").Value = "'" + $tarih
    $ws.Range("C$row").Value = "'" + $birim
    $ws.Range("D$row").Value = "'" + $parselSayi
    $ws.Range("E$row").Value = "'" + $is
    $ws.Range("F$row").Value = "'" + $personel

    # Drop the implicit "quote prefix" style so the new cells keep the
    # sheet's plain (unstyled) formatting, like every other data row.
    $ws.Range("A$row`:F$row").Style = "Normal"
}
